$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New recalculated NATMI values (per Dr Hou's advice) for rows 2-10.
# Columns: E,G,H,I,J,K,M,N,O,P,Q,R,S,T  (F and L stay unchanged at 1)
$data = @{
    2  = @{ E=3; G=142.9073533333333;  H=428.72206;         I=0.5576664151504187; J=0.5576664151504188; K=3; M=1.104136666666667;  N=3.31241;           O=0.2772362398998524; P=0.2772362398998525; Q=157.7892487516222;  R=1420.1032387646;   S=0.1546053400547322;  T=0.1546053400547322 }
    3  = @{ E=3; G=142.9073533333333;  H=428.72206;         I=0.5576664151504187; J=0.5576664151504188; K=3; M=2.230986666666666;  N=6.692959999999999; O=0.5601755411317187; P=0.5601755411317187; Q=318.8243998552888;  R=2869.4195986976;   S=0.3123910858778715;  T=0.3123910858778716 }
    4  = @{ E=3; G=142.9073533333333;  H=428.72206;         I=0.5576664151504187; J=0.5576664151504188; K=3; M=0.647533;           N=1.942599;           O=0.1625882189684289; P=0.1625882189684289; Q=92.53722722599333;  R=832.83504503394;   S=0.09066998921781504; T=0.09066998921781505 }
    5  = @{ E=3; G=63.967809;           H=191.903427;         I=0.2496211559306514; J=0.2496211559306514; K=3; M=1.104136666666667;  N=3.31241;           O=0.2772362398998524; P=0.2772362398998525; Q=70.62920340322999;  R=635.66283062907;   S=0.06920403066966853; T=0.06920403066966856 }
    6  = @{ E=3; G=63.967809;           H=191.903427;         I=0.2496211559306514; J=0.2496211559306514; K=3; M=2.230986666666666;  N=6.692959999999999; O=0.5601755411317187; P=0.5601755411317187; Q=142.71132897488;    R=1284.40196077392;  S=0.1398316661013778;  T=0.1398316661013778 }
    7  = @{ E=3; G=63.967809;           H=191.903427;         I=0.2496211559306514; J=0.2496211559306514; K=3; M=0.647533;           N=1.942599;           O=0.1625882189684289; P=0.1625882189684289; Q=41.421267265197;    R=372.791405386773;  S=0.04058545915960507; T=0.04058545915960508 }
    8  = @{ E=3; G=49.38440333333333;  H=148.15321;           I=0.1927124289189298; J=0.1927124289189298; K=3; M=1.104136666666667;  N=3.31241;           O=0.2772362398998524; P=0.2772362398998525; Q=54.52713048178888;  R=490.7441743361;    S=0.05342686917545168; T=0.0534268691754517 }
    9  = @{ E=3; G=49.38440333333333;  H=148.15321;           I=0.1927124289189298; J=0.1927124289189298; K=3; M=2.230986666666666;  N=6.692959999999999; O=0.5601755411317187; P=0.5601755411317187; Q=110.1759453779555;  R=991.5835084015999; S=0.1079527891524694;  T=0.1079527891524694 }
    10 = @{ E=3; G=49.38440333333333;  H=148.15321;           I=0.1927124289189298; J=0.1927124289189298; K=3; M=0.647533;           N=1.942599;           O=0.1625882189684289; P=0.1625882189684289; Q=31.97803084364333;  R=287.80227759279;   S=0.03133277059100874; T=0.03133277059100875 }
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    $ws.Range("E$r").Value = $row.E
    $ws.Range("G$r").Value = $row.G
    $ws.Range("H$r").Value = $row.H
    $ws.Range("I$r").Value = $row.I
    $ws.Range("J$r").Value = $row.J
    $ws.Range("K$r").Value = $row.K
    $ws.Range("M$r").Value = $row.M
    $ws.Range("N$r").Value = $row.N
    $ws.Range("O$r").Value = $row.O
    $ws.Range("P$r").Value = $row.P
    $ws.Range("Q$r").Value = $row.Q
    $ws.Range("R$r").Value = $row.R
    $ws.Range("S$r").Value = $row.S
    $ws.Range("T$r").Value = $row.T
}
